$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.714211225509644
$ws.Range("C1").Value = 2.008512496948242
$ws.Range("D1").Value = 1.573822736740112
$ws.Range("E1").Value = 1.430870890617371
